# Updates cryptos list values (Price / Volume(1h) columns) per the
# "Updated cryptos list" GitHub Actions commit.
#
# These cells are stored as literal text (t="inlineStr" in the OOXML),
# even though many of the "Price" values look like plain numbers
# (e.g. "324.43", "0.4563"). Assigning such a string straight to
# .Value would make Excel auto-convert it into a real number, which
# would change the cell's stored type and introduce floating-point
# artifacts. To keep these as genuine text values we temporarily force
# the cell's number format to Text ("@") before assigning the value,
# then clear the format again afterwards so no extra cell styling is
# left behind (matching the original, style-less cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '28.892.28' },
    @{ Cell = 'E2'; Value = '  -0.27%  ' },
    @{ Cell = 'D3'; Value = '1.916.69' },
    @{ Cell = 'E3'; Value = '  +0.67%  ' },
    @{ Cell = 'E4'; Value = '  -0.03%  ' },
    @{ Cell = 'D5'; Value = '324.43' },
    @{ Cell = 'E5'; Value = '  -0.06%  ' },
    @{ Cell = 'E6'; Value = '  +0.01%  ' },
    @{ Cell = 'D7'; Value = '0.4563' },
    @{ Cell = 'E8'; Value = '  -0.46%  ' },
    @{ Cell = 'D9'; Value = '0.07739' },
    @{ Cell = 'E9'; Value = '  +0.43%  ' },
    @{ Cell = 'D10'; Value = '0.9753' },
    @{ Cell = 'E10'; Value = '  -0.54%  ' },
    @{ Cell = 'D11'; Value = '22.28' },
    @{ Cell = 'E11'; Value = '  +1.18%  ' },
    @{ Cell = 'D12'; Value = '1.922.78' },
    @{ Cell = 'E12'; Value = '  -0.01%  ' },
    @{ Cell = 'D13'; Value = '5.690' },
    @{ Cell = 'E13'; Value = '  +0.35%  ' },
    @{ Cell = 'D14'; Value = '6.962' },
    @{ Cell = 'D15'; Value = '0.07001' },
    @{ Cell = 'E15'; Value = '  -0.54%  ' },
    @{ Cell = 'D16'; Value = '1.005' },
    @{ Cell = 'E16'; Value = '  +0.03%  ' },
    @{ Cell = 'D17'; Value = '84.28' },
    @{ Cell = 'E17'; Value = '  +0.52%  ' },
    @{ Cell = 'D18'; Value = '0.000009468' },
    @{ Cell = 'E18'; Value = '  -0.64%  ' },
    @{ Cell = 'E19'; Value = '  -0.52%  ' },
    @{ Cell = 'E20'; Value = '  -0.03%  ' },
    @{ Cell = 'D21'; Value = '28.913.51' },
    @{ Cell = 'E21'; Value = '  -0.16%  ' },
    @{ Cell = 'D22'; Value = '5.337' },
    @{ Cell = 'E22'; Value = '  +0.15%  ' },
    @{ Cell = 'E23'; Value = '  +1.51%  ' },
    @{ Cell = 'D24'; Value = '2.171.12' },
    @{ Cell = 'E24'; Value = '  +0.83%  ' },
    @{ Cell = 'E25'; Value = '  -2.22%  ' },
    @{ Cell = 'D26'; Value = '157.61' },
    @{ Cell = 'E26'; Value = '  +0.31%  ' },
    @{ Cell = 'D27'; Value = '19.04' },
    @{ Cell = 'E27'; Value = '  -0.44%  ' },
    @{ Cell = 'D28'; Value = '5.612' },
    @{ Cell = 'E28'; Value = '  +0.41%  ' },
    @{ Cell = 'D29'; Value = '117.81' },
    @{ Cell = 'E29'; Value = '  +0.16%  ' },
    @{ Cell = 'D30'; Value = '1.840' },
    @{ Cell = 'E30'; Value = '  -0.62%  ' },
    @{ Cell = 'D31'; Value = '0.09306' },
    @{ Cell = 'E31'; Value = '  +0.37%  ' },
    @{ Cell = 'D32'; Value = '0.8686' },
    @{ Cell = 'E32'; Value = '  +0.99%  ' },
    @{ Cell = 'D33'; Value = '5.092' },
    @{ Cell = 'E33'; Value = '  +0.27%  ' },
    @{ Cell = 'D34'; Value = '1.240' },
    @{ Cell = 'E34'; Value = '  -0.68%  ' },
    @{ Cell = 'D35'; Value = '3.008' },
    @{ Cell = 'E35'; Value = '  +0.40%  ' },
    @{ Cell = 'D36'; Value = '0.05687' },
    @{ Cell = 'E36'; Value = '  -0.03%  ' },
    @{ Cell = 'D37'; Value = '1.147' },
    @{ Cell = 'E37'; Value = '  -0.03%  ' },
    @{ Cell = 'E38'; Value = '  -0.09%  ' },
    @{ Cell = 'D39'; Value = '0.02037' },
    @{ Cell = 'E39'; Value = '  +0.12%  ' },
    @{ Cell = 'D40'; Value = '3.056' },
    @{ Cell = 'E40'; Value = '  +10.68%  ' },
    @{ Cell = 'D41'; Value = '7.478' },
    @{ Cell = 'E41'; Value = '  +0.19%  ' },
    @{ Cell = 'D42'; Value = '0.5488' },
    @{ Cell = 'E42'; Value = '  -0.41%  ' },
    @{ Cell = 'D43'; Value = '0.1754' },
    @{ Cell = 'E43'; Value = '  -0.04%  ' },
    @{ Cell = 'D44'; Value = '9.322' },
    @{ Cell = 'E44'; Value = '  +0.65%  ' },
    @{ Cell = 'D45'; Value = '0.000002845' },
    @{ Cell = 'E45'; Value = '  +17.37%  ' },
    @{ Cell = 'D46'; Value = '2.156' },
    @{ Cell = 'E46'; Value = '  +3.44%  ' },
    @{ Cell = 'D47'; Value = '0.5162' },
    @{ Cell = 'E47'; Value = '  -0.55%  ' },
    @{ Cell = 'D48'; Value = '0.06927' },
    @{ Cell = 'E48'; Value = '  +1.53%  ' },
    @{ Cell = 'D49'; Value = '11.07' },
    @{ Cell = 'E49'; Value = '  -2.41%  ' },
    @{ Cell = 'D50'; Value = '110.65' },
    @{ Cell = 'E50'; Value = '  -0.57%  ' },
    @{ Cell = 'D51'; Value = '1.761' },
    @{ Cell = 'E51'; Value = '  -0.86%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.ClearFormats()
}
